$d = $word.ActiveDocument

function Split-RunBoundary($rangeObj) {
    # Toggling Bold on/off forces the interop engine to materialize a
    # distinct run for exactly this sub-range without altering the
    # visible formatting (Bold ends up False again, matching the rest
    # of the paragraph).
    $rangeObj.Font.Bold = 1
    $rangeObj.Font.Bold = 0
}

# ---------------------------------------------------------------------
# Change 1: "... Sistema salva aluguel e exibe página com as informações
# de cobrança" -> "... Sistema salva aluguel e exibe uma página de
# impressão com as informações de cobrança"
# The old run "página com as informações de cobrança" is split into four
# runs: "uma " | "página" | " de impressão" | " com as informações de
# cobrança" (all identical formatting, just separate <w:r> elements).
# ---------------------------------------------------------------------
$f1 = $d.Content
$f1.Find.Execute("página com as informações de cobrança", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($f1.Find.Found) {
    $start = $f1.Start
    $end = $f1.End

    # Insert "uma " before "página"
    $whole = $d.Range($start, $end)
    $whole.InsertBefore("uma ")

    # Insert " de impressão" right after "página" (6 chars), before
    # " com as informações de cobrança".
    $paginaEnd = $start + ("uma ").Length + ("página").Length
    $tail = $d.Range($paginaEnd, $end + ("uma ").Length)
    $tail.InsertBefore(" de impressão")

    # Recompute the four final boundaries and force each into its own run.
    $b0 = $start
    $b1 = $b0 + ("uma ").Length
    $b2 = $b1 + ("página").Length
    $b3 = $b2 + (" de impressão").Length
    $b4 = $b3 + (" com as informações de cobrança").Length

    Split-RunBoundary ($d.Range($b0, $b1))
    Split-RunBoundary ($d.Range($b1, $b2))
    Split-RunBoundary ($d.Range($b2, $b3))
    Split-RunBoundary ($d.Range($b3, $b4))
}

# ---------------------------------------------------------------------
# Change 2: "Visitante escolhe gerar pdf para imprimir informações do
# aluguel e pagamento" -> "Visitante escolhe imprimir informações do
# aluguel e pagamento"
# "Visitante escolhe " + "gerar pdf para imprimir" collapse down to
# "Visitante escolhe" | " imprimir" (two runs); the trailing
# " informações do aluguel e pagamento" run is untouched.
# ---------------------------------------------------------------------
$f2 = $d.Content
$f2.Find.Execute("Visitante escolhe gerar pdf para imprimir", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($f2.Find.Found) {
    $start = $f2.Start
    $end = $f2.End

    $whole = $d.Range($start, $end)
    $whole.Text = "Visitante escolhe imprimir"

    $b0 = $start
    $b1 = $b0 + ("Visitante escolhe").Length
    $b2 = $b1 + (" imprimir").Length

    Split-RunBoundary ($d.Range($b0, $b1))
    Split-RunBoundary ($d.Range($b1, $b2))
}

# ---------------------------------------------------------------------
# Change 3: "Sistema redireciona o visitante para uma página de
# impressão do pdf" -> "Sistema retorna que o aluguel foi efetuado com
# sucesso"
# The single old run splits into "Sistema " | "retorna que o aluguel
# foi efetuado com sucesso".
# ---------------------------------------------------------------------
$f3 = $d.Content
$f3.Find.Execute("Sistema redireciona o visitante para uma página de impressão do pdf", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($f3.Find.Found) {
    $start = $f3.Start
    $end = $f3.End

    $whole = $d.Range($start, $end)
    $whole.Text = "Sistema retorna que o aluguel foi efetuado com sucesso"

    $b0 = $start
    $b1 = $b0 + ("Sistema ").Length
    $b2 = $b1 + ("retorna que o aluguel foi efetuado com sucesso").Length

    Split-RunBoundary ($d.Range($b0, $b1))
    Split-RunBoundary ($d.Range($b1, $b2))
}
